$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the default (unstyled) cell style, used to strip
# the quote-prefix style Excel applies when a value is entered with a
# leading apostrophe (forcing text interpretation of numeric-looking strings).
$normalStyle = $ws.Cells.Item(1, 1).Style

$ws.Cells.Item(2, 4).Value = "42.203.82"
$ws.Cells.Item(2, 5).Value = "  +1.05%  "
$ws.Cells.Item(3, 4).Value = "2.173.56"
$ws.Cells.Item(3, 5).Value = "  -0.07%  "
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).Value = "'252.89"
$ws.Cells.Item(5, 4).Style = $normalStyle
$ws.Cells.Item(5, 5).Value = "  +6.03%  "
$ws.Cells.Item(6, 4).Value = "'0.610"
$ws.Cells.Item(6, 4).Style = $normalStyle
$ws.Cells.Item(6, 5).Value = "  -0.24%  "
$ws.Cells.Item(7, 5).Value = "  +0.63%  "
$ws.Cells.Item(8, 5).Value = "  -0.02%  "
$ws.Cells.Item(9, 5).Value = "  -0.35%  "
$ws.Cells.Item(10, 4).Value = "'40.47"
$ws.Cells.Item(10, 4).Style = $normalStyle
$ws.Cells.Item(10, 5).Value = "  -0.41%  "
$ws.Cells.Item(11, 4).Value = "'0.0909"
$ws.Cells.Item(11, 4).Style = $normalStyle
$ws.Cells.Item(11, 5).Value = "  -0.35%  "
$ws.Cells.Item(12, 4).Value = "'6.78"
$ws.Cells.Item(12, 4).Style = $normalStyle
$ws.Cells.Item(12, 5).Value = "  +0.45%  "
$ws.Cells.Item(13, 5).Value = "  +0.57%  "
$ws.Cells.Item(14, 4).Value = "2.499.24"
$ws.Cells.Item(14, 5).Value = "  -0.09%  "
$ws.Cells.Item(15, 4).Value = "'14.17"
$ws.Cells.Item(15, 4).Style = $normalStyle
$ws.Cells.Item(15, 5).Value = "  -1.40%  "
$ws.Cells.Item(16, 4).Value = "2.164.81"
$ws.Cells.Item(16, 5).Value = "  -0.44%  "
$ws.Cells.Item(17, 4).Value = "'0.763"
$ws.Cells.Item(17, 4).Style = $normalStyle
$ws.Cells.Item(17, 5).Value = "  -2.68%  "
$ws.Cells.Item(18, 4).Value = "42.095.14"
$ws.Cells.Item(18, 5).Value = "  +1.10%  "
$ws.Cells.Item(19, 5).Value = "  -0.75%  "
$ws.Cells.Item(20, 4).Value = "'70.60"
$ws.Cells.Item(20, 4).Style = $normalStyle
$ws.Cells.Item(20, 5).Value = "  +0.62%  "
$ws.Cells.Item(21, 5).Value = "  +0.49%  "
$ws.Cells.Item(22, 4).Value = "'226.35"
$ws.Cells.Item(22, 4).Style = $normalStyle
$ws.Cells.Item(22, 5).Value = "  -0.22%  "
$ws.Cells.Item(23, 5).Value = "  -5.07%  "
$ws.Cells.Item(24, 5).Value = "  +5.50%  "
$ws.Cells.Item(25, 5).Value = "  -0.25%  "
$ws.Cells.Item(26, 5).Value = "  -3.43%  "
$ws.Cells.Item(27, 4).Value = "'3.31"
$ws.Cells.Item(27, 4).Style = $normalStyle
$ws.Cells.Item(27, 5).Value = "  +1.64%  "
$ws.Cells.Item(28, 2).Value = "Toncoin"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(28, 4).Value = "'2.22"
$ws.Cells.Item(28, 4).Style = $normalStyle
$ws.Cells.Item(28, 5).Value = "  +1.23%  "
$ws.Cells.Item(29, 2).Value = "PancakeSwap"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(29, 4).Value = "'2.19"
$ws.Cells.Item(29, 4).Style = $normalStyle
$ws.Cells.Item(29, 5).Value = "  -0.55%  "
$ws.Cells.Item(30, 4).Value = "'36.66"
$ws.Cells.Item(30, 4).Style = $normalStyle
$ws.Cells.Item(30, 5).Value = "  +11.84%  "
$ws.Cells.Item(31, 4).Value = "'168.99"
$ws.Cells.Item(31, 4).Style = $normalStyle
$ws.Cells.Item(31, 5).Value = "  -1.32%  "
$ws.Cells.Item(32, 4).Value = "'19.99"
$ws.Cells.Item(32, 4).Style = $normalStyle
$ws.Cells.Item(32, 5).Value = "  +0.53%  "
$ws.Cells.Item(33, 4).Value = "'0.0803"
$ws.Cells.Item(33, 4).Style = $normalStyle
$ws.Cells.Item(33, 5).Value = "  +3.22%  "
$ws.Cells.Item(34, 4).Value = "'5.12"
$ws.Cells.Item(34, 4).Style = $normalStyle
$ws.Cells.Item(34, 5).Value = "  -3.40%  "
$ws.Cells.Item(35, 5).Value = "  -0.34%  "
$ws.Cells.Item(36, 5).Value = "  +3.96%  "
$ws.Cells.Item(37, 5).Value = "  -2.52%  "
$ws.Cells.Item(38, 4).Value = "'0.0329"
$ws.Cells.Item(38, 4).Style = $normalStyle
$ws.Cells.Item(38, 5).Value = "  +5.59%  "
$ws.Cells.Item(39, 4).Value = "'11.95"
$ws.Cells.Item(39, 4).Style = $normalStyle
$ws.Cells.Item(39, 5).Value = "  -2.18%  "
$ws.Cells.Item(40, 5).Value = "  -2.81%  "
$ws.Cells.Item(41, 4).Value = "'0.195"
$ws.Cells.Item(41, 4).Style = $normalStyle
$ws.Cells.Item(41, 5).Value = "  +3.02%  "
$ws.Cells.Item(42, 4).Value = "'58.87"
$ws.Cells.Item(42, 4).Style = $normalStyle
$ws.Cells.Item(42, 5).Value = "  -0.62%  "
$ws.Cells.Item(43, 5).Value = "  -5.36%  "
$ws.Cells.Item(44, 4).Value = "'102.17"
$ws.Cells.Item(44, 4).Style = $normalStyle
$ws.Cells.Item(44, 5).Value = "  +4.36%  "
$ws.Cells.Item(45, 2).Value = "WOONetwork"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Cells.Item(45, 4).Value = "'0.467"
$ws.Cells.Item(45, 4).Style = $normalStyle
$ws.Cells.Item(45, 5).Value = "  +14.27%  "
$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(46, 4).Value = "'8.28"
$ws.Cells.Item(46, 4).Style = $normalStyle
$ws.Cells.Item(46, 5).Value = "  -2.41%  "
$ws.Cells.Item(47, 2).Value = "NEARProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(47, 4).Value = "'2.44"
$ws.Cells.Item(47, 4).Style = $normalStyle
$ws.Cells.Item(47, 5).Value = "  +10.58%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).Value = "'0.0968"
$ws.Cells.Item(48, 4).Style = $normalStyle
$ws.Cells.Item(48, 5).Value = "  +0.14%  "
$ws.Cells.Item(49, 4).Value = "'1.08"
$ws.Cells.Item(49, 4).Style = $normalStyle
$ws.Cells.Item(49, 5).Value = "  +0.11%  "
$ws.Cells.Item(50, 5).Value = "  +0.49%  "
$ws.Cells.Item(51, 5).Value = "  +0.80%  "
